# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.335.80'
$ws.Range('E2').Value = '  -2.99%  '
$ws.Range('D3').Value = '3.619.12'
$ws.Range('E3').Value = '  -3.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '587.74'
$ws.Range('E5').Value = '  -2.35%  '
$ws.Range('D6').Value = '181.59'
$ws.Range('E6').Value = '  -2.35%  '
$ws.Range('D7').Value = '0.611'
$ws.Range('E7').Value = '  -3.77%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '0.673'
$ws.Range('E9').Value = '  -7.28%  '
$ws.Range('D10').Value = '0.145'
$ws.Range('E10').Value = '  -11.25%  '
$ws.Range('D11').Value = '53.34'
$ws.Range('E11').Value = '  -6.45%  '
$ws.Range('D12').Value = '0.0000253'
$ws.Range('E12').Value = '  -14.26%  '
$ws.Range('D13').Value = '9.93'
$ws.Range('E13').Value = '  -8.20%  '
$ws.Range('D14').Value = '4.208.37'
$ws.Range('E14').Value = '  -3.19%  '
$ws.Range('D15').Value = '3.629.49'
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('D17').Value = '67.250.33'
$ws.Range('E17').Value = '  -2.95%  '
$ws.Range('D18').Value = '18.34'
$ws.Range('E18').Value = '  -6.13%  '
$ws.Range('D19').Value = '12.23'
$ws.Range('E19').Value = '  -5.70%  '
$ws.Range('D20').Value = '1.06'
$ws.Range('E20').Value = '  -5.95%  '
$ws.Range('D21').Value = '395.97'
$ws.Range('E21').Value = '  -4.35%  '
$ws.Range('D22').Value = '4.31'
$ws.Range('E22').Value = '  -7.45%  '
$ws.Range('D23').Value = '85.35'
$ws.Range('E23').Value = '  -4.70%  '
$ws.Range('D24').Value = '2.85'
$ws.Range('E24').Value = '  -6.92%  '
$ws.Range('D25').Value = '12.27'
$ws.Range('E25').Value = '  -4.97%  '
$ws.Range('D26').Value = '6.06'
$ws.Range('E26').Value = '  -0.42%  '
$ws.Range('D27').Value = '10.26'
$ws.Range('E27').Value = '  -7.04%  '
$ws.Range('D28').Value = '3.61'
$ws.Range('E28').Value = '  -8.83%  '
$ws.Range('D29').Value = '8.98'
$ws.Range('E29').Value = '  -5.82%  '
$ws.Range('D30').Value = '31.15'
$ws.Range('E30').Value = '  -5.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.70'
$ws.Range('E31').Value = '  -9.13%  '
$ws.Range('D32').Value = '66.52'
$ws.Range('E32').Value = '  +1.79%  '
$ws.Range('D33').Value = '11.87'
$ws.Range('E33').Value = '  -5.39%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.112'
$ws.Range('E34').Value = '  -5.73%  '
$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').Value = '589.07'
$ws.Range('E35').Value = '  -4.14%  '
$ws.Range('D36').Value = '41.45'
$ws.Range('E36').Value = '  -6.07%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('D39').Value = '0.375'
$ws.Range('E39').Value = '  -7.75%  '
$ws.Range('D40').Value = '0.0₃0735'
$ws.Range('E40').Value = '  -18.80%  '
$ws.Range('D41').Value = '0.132'
$ws.Range('E41').Value = '  -4.01%  '
$ws.Range('D42').Value = '2.77'
$ws.Range('E42').Value = '  -9.87%  '
$ws.Range('D43').Value = '0.0411'
$ws.Range('E43').Value = '  -7.69%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.700.51'
$ws.Range('E44').Value = '  -3.25%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').Value = '2.38'
$ws.Range('E45').Value = '  -13.39%  '
$ws.Range('E46').Value = '  -4.50%  '
$ws.Range('D47').Value = '3.07'
$ws.Range('E47').Value = '  -3.81%  '
$ws.Range('E48').Value = '  -6.70%  '
$ws.Range('D49').Value = '137.13'
$ws.Range('E49').Value = '  -3.53%  '
$ws.Range('D50').Value = '8.19'
$ws.Range('E50').Value = '  -12.58%  '
$ws.Range('D51').Value = '2.58'
$ws.Range('E51').Value = '  -7.45%  '
